{"js": "// The document has a table cell that contains two paragraphs:\n//   \"Click the Choose New Sprite icon from the Sprites bar.\"\n//   \"Select Mouse1\"\n// The second paragraph is a leftover/incorrect instruction (telling the\n// reader to pick the built-in \"Mouse1\" sprite) and must be removed\n// entirely, per the commit message \"remove mistake instruction to add\n// wrong sprite\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"Select Mouse1\";\nconst toRemove = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === target) {\n    toRemove.push(paragraphs.items[i]);\n  }\n}\n\nfor (const para of toRemove) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# The document has a table cell containing two paragraphs:\n#   \"Click the Choose New Sprite icon from the Sprites bar.\"\n#   \"Select Mouse1\"\n# The second paragraph is a leftover/incorrect instruction (telling the\n# reader to pick the built-in \"Mouse1\" sprite) and must be removed\n# entirely, per the commit message \"remove mistake instruction to add\n# wrong sprite\".\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is (or contains) \"Select Mouse1\".\n$count = $d.Paragraphs.Count\n$targetIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Select Mouse1*\") {\n        $targetIdx = $i\n    }\n}\n\nif ($targetIdx -gt 0) {\n    $p = $d.Paragraphs.Item($targetIdx)\n    $r = $p.Range\n    $len = $r.End - $r.Start\n    # Delete the whole paragraph, including its trailing paragraph mark,\n    # so the <w:p> element itself is removed rather than left empty.\n    $r.Delete(1, $len)\n}\n"}
